$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 47828.5
$ws.Range("J3").Value = 47828.5
$ws.Range("L3").Value = 47828.5
$ws.Range("N3").Value = -48056.5
$ws.Range("H10").Value = 20100
$ws.Range("I10").Value = 300
$ws.Range("K10").Value = 300
$ws.Range("M10").Value = -7
$ws.Range("H19").Value = 330.5
$ws.Range("J19").Value = 318.27274
$ws.Range("L19").Value = 318.27274
$ws.Range("N19").Value = -668.27274
$ws.Range("H102").Value = 47828.5
$ws.Range("J102").Value = 47828.5
$ws.Range("L102").Value = 47828.5
$ws.Range("N102").Value = -54318.5
$ws.Range("H113").Value = 8740.823
$ws.Range("I113").Value = 8450
$ws.Range("J113").Value = 8830.308000000001
$ws.Range("K113").Value = 8450
$ws.Range("L113").Value = 8830.308000000001
$ws.Range("M113").Value = -5196
$ws.Range("N113").Value = -15338.308
$ws.Range("H132").Value = 1768.2424
$ws.Range("I132").Value = 1865.2
$ws.Range("J132").Value = 798.6667
$ws.Range("K132").Value = 5595.6
$ws.Range("L132").Value = 2396.0001
$ws.Range("M132").Value = -3065.6
$ws.Range("N132").Value = -7456.0001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6608.0166
$ws.Range("I32").Value = 6724.381
$ws.Range("J32").Value = 6350.7896
$ws.Range("K32").Value = 6724.381
$ws.Range("L32").Value = 6350.7896
$ws.Range("M32").Value = -6437.381
$ws.Range("N32").Value = -6924.7896
$ws.Range("H106").Value = 51672
$ws.Range("J106").Value = 51672
$ws.Range("L106").Value = 51672
$ws.Range("N106").Value = -54196
$ws.Range("H132").Value = 2125.9033
$ws.Range("I132").Value = 1225.3024
$ws.Range("J132").Value = 4164.1055
$ws.Range("K132").Value = 3675.9072
$ws.Range("L132").Value = 12492.3165
$ws.Range("M132").Value = -1145.9072
$ws.Range("N132").Value = -17552.3165
$ws.Range("H140").Value = 46500
$ws.Range("J140").Value = 46500
$ws.Range("L140").Value = 46500
$ws.Range("N140").Value = -56860

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H140").Value = 41352.766
$ws.Range("J140").Value = 41352.766
$ws.Range("L140").Value = 41352.766
$ws.Range("N140").Value = -51712.766

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H95").Value = 8513.714
$ws.Range("J95").Value = 8513.714
$ws.Range("L95").Value = 8513.714
$ws.Range("N95").Value = -14005.714

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 2476.7385
$ws.Range("I68").Value = 3202.1191
$ws.Range("J68").Value = 1814.4348
$ws.Range("K68").Value = 9606.3573
$ws.Range("L68").Value = 5443.3044
$ws.Range("M68").Value = -8795.3573
$ws.Range("N68").Value = -7065.3044
$ws.Range("H71").Value = 2476.7385
$ws.Range("I71").Value = 3202.1191
$ws.Range("J71").Value = 1814.4348
$ws.Range("K71").Value = 28819.0719
$ws.Range("L71").Value = 16329.9132
$ws.Range("M71").Value = -24763.0719
$ws.Range("N71").Value = -24441.9132
$ws.Range("H80").Value = 4555.6665
$ws.Range("I80").Value = 2425.25
$ws.Range("J80").Value = 6260
$ws.Range("K80").Value = 7275.75
$ws.Range("L80").Value = 18780
$ws.Range("M80").Value = -6339.75
$ws.Range("N80").Value = -20652
$ws.Range("H83").Value = 4555.6665
$ws.Range("I83").Value = 2425.25
$ws.Range("J83").Value = 6260
$ws.Range("K83").Value = 21827.25
$ws.Range("L83").Value = 56340
$ws.Range("M83").Value = -17147.25
$ws.Range("N83").Value = -65700
$ws.Range("H137").Value = 41686716
$ws.Range("J137").Value = 50023380
$ws.Range("L137").Value = 150070140
$ws.Range("N137").Value = -150080340

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H92").Value = 9999
$ws.Range("J92").Value = 9999
$ws.Range("L92").Value = 9999
$ws.Range("N92").Value = -13743
$ws.Range("H95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("L95").Value = 0
$ws.Range("H132").Value = 3516.8125
$ws.Range("I132").Value = 3379.4092
$ws.Range("J132").Value = 3633.077
$ws.Range("K132").Value = 10138.2276
$ws.Range("L132").Value = 10899.231
$ws.Range("M132").Value = -7608.2276
$ws.Range("N132").Value = -15959.231

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H106").Value = 13602.2
$ws.Range("J106").Value = 13602.2
$ws.Range("L106").Value = 13602.2
$ws.Range("N106").Value = -16126.2
$ws.Range("H122").Value = 16294714
$ws.Range("I122").Value = 23817856
$ws.Range("J122").Value = 5010000
$ws.Range("K122").Value = 71453568
$ws.Range("L122").Value = 15030000
$ws.Range("M122").Value = -71451118
$ws.Range("N122").Value = -15034900
$ws.Range("H136").Value = 6630.9805
$ws.Range("I136").Value = 4927.6284
$ws.Range("J136").Value = 10357.0625
$ws.Range("K136").Value = 14782.8852
$ws.Range("L136").Value = 31071.1875
$ws.Range("M136").Value = -12232.8852
$ws.Range("N136").Value = -36171.1875

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2001.2222
$ws.Range("I81").Value = 1702.2
$ws.Range("J81").Value = 2375
$ws.Range("K81").Value = 3404.4
$ws.Range("L81").Value = 4750
$ws.Range("M81").Value = -2343.4
$ws.Range("N81").Value = -6872
$ws.Range("H84").Value = 2001.2222
$ws.Range("I84").Value = 1702.2
$ws.Range("J84").Value = 2375
$ws.Range("K84").Value = 17022
$ws.Range("L84").Value = 23750
$ws.Range("M84").Value = -11718
$ws.Range("N84").Value = -34358
$ws.Range("H101").Value = 16700.4
$ws.Range("J101").Value = 16700.4
$ws.Range("L101").Value = 16700.4
$ws.Range("N101").Value = -23190.4
$ws.Range("H105").Value = 48666.25
$ws.Range("J105").Value = 48666.25
$ws.Range("L105").Value = 48666.25
$ws.Range("N105").Value = -55654.25
$ws.Range("H122").Value = 15501
$ws.Range("I122").Value = 14001.333
$ws.Range("J122").Value = 20000
$ws.Range("K122").Value = 42003.999
$ws.Range("L122").Value = 60000
$ws.Range("M122").Value = -39553.999
$ws.Range("N122").Value = -64900
$ws.Range("H126").Value = 827.96295
$ws.Range("I126").Value = 623.73914
$ws.Range("J126").Value = 2002.25
$ws.Range("K126").Value = 1871.21742
$ws.Range("L126").Value = 6006.75
$ws.Range("M126").Value = 598.7825800000001
$ws.Range("N126").Value = -10946.75
$ws.Range("H136").Value = 1563.7778
$ws.Range("I136").Value = 941.5333000000001
$ws.Range("J136").Value = 2341.5833
$ws.Range("K136").Value = 2824.5999
$ws.Range("L136").Value = 7024.749899999999
$ws.Range("M136").Value = -274.5999000000002
$ws.Range("N136").Value = -12124.7499
$ws.Range("H138").Value = 24995.666
$ws.Range("J138").Value = 24995.666
$ws.Range("L138").Value = 24995.666
$ws.Range("N138").Value = -35275.666

# Cell deletions (cells removed entirely in the target revision)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("N95").ClearContents()
